# Edit: Mon, Jun 01, 2020 11:04:58 PM
#
# 1) Slide 16's table (the "Total Outflow" plenary table) is switched from
#    the deck's custom "Table_0" style to PowerPoint's built-in
#    "No Style, Table Grid" style.
# 2) The presentation's design is switched from the "Integral" theme to the
#    built-in "Office Theme" - i.e. the slide master's theme (theme1.xml)
#    picks up the Office Theme's 12 scheme colors.

$p = $ppt.ActivePresentation

# --- 1) Table style swap on slide 16, shape 3 (the graphicFrame/table) ---
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{97C58DA9-4482-46CA-8809-069433AE8323}")

# --- 2) Switch the deck's theme colours from Integral to Office Theme ---
# ThemeColorScheme item order is fixed: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink - matches <a:clrScheme> child order in theme1.xml.
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
